$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price values in column D must be kept as text
# (they are stored as strings in the source data), so force text
# number format before assigning to avoid float conversion / precision loss.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "243.68"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "24.04"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.05753"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "6.500"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.8157"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.8507"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.1353"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.06957"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.03151"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.02872"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.756"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.001516"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.04684"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.0006016"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.006283"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "0.001238"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "0.004292"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.00008728"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.3171"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.1339"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.0002332"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.03628"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.1052"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.002809"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.006331"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.007534"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.00005271"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.2903"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.002343"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.00002102"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.0002002"

# Text column updates (Coin name, Link, Volume label) - row reorder for
# KickToken / BKEXToken / CEJI plus small label refreshes.
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICK"
